# Corrected financial figures for LG유플러스 IFRS list (error solve ifrs list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2014/12  (IFRS연결)
$ws.Range("D2").Value = 109998
$ws.Range("E2").Value = 5763
$ws.Range("F2").Value = 5763
$ws.Range("G2").Value = 3201
$ws.Range("H2").Value = 2277
$ws.Range("I2").Value = 2282
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 120127
$ws.Range("L2").Value = 78348
$ws.Range("M2").Value = 41778
$ws.Range("N2").Value = 41771
$ws.Range("O2").Value = 7
$ws.Range("P2").Value = 25740
$ws.Range("Q2").Value = 20154
$ws.Range("R2").Value = -23073
$ws.Range("S2").Value = 3095
$ws.Range("T2").Value = 21448
$ws.Range("U2").Value = -1294
$ws.Range("V2").Value = 49157
$ws.Range("W2").Value = 5.24
$ws.Range("X2").Value = 2.07
$ws.Range("Y2").Value = 5.57
$ws.Range("Z2").Value = 1.92
$ws.Range("AA2").Value = 187.53
$ws.Range("AB2").Value = 62.22
$ws.Range("AC2").Value = 523
$ws.Range("AD2").Value = 22
$ws.Range("AE2").Value = 9567
$ws.Range("AF2").Value = 1.2
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.3
$ws.Range("AI2").Value = 28.7
$ws.Range("AJ2").Value = 436611361

# Row 3: 2015/12  (IFRS연결)
$ws.Range("D3").Value = 107952
$ws.Range("E3").Value = 6323
$ws.Range("F3").Value = 6323
$ws.Range("G3").Value = 4659
$ws.Range("H3").Value = 3512
$ws.Range("I3").Value = 3514
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 119510
$ws.Range("L3").Value = 75026
$ws.Range("M3").Value = 44484
$ws.Range("N3").Value = 44480
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 25740
$ws.Range("Q3").Value = 17929
$ws.Range("R3").Value = -15107
$ws.Range("S3").Value = -4060
$ws.Range("T3").Value = 13754
$ws.Range("U3").Value = 4175
$ws.Range("V3").Value = 45752
$ws.Range("W3").Value = 5.86
$ws.Range("X3").Value = 3.25
$ws.Range("Y3").Value = 8.15
$ws.Range("Z3").Value = 2.93
$ws.Range("AA3").Value = 168.66
$ws.Range("AB3").Value = 72.76000000000001
$ws.Range("AC3").Value = 805
$ws.Range("AD3").Value = 12.92
$ws.Range("AE3").Value = 10187
$ws.Range("AF3").Value = 1.02
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 2.4
$ws.Range("AI3").Value = 31.06
$ws.Range("AJ3").Value = 436611361

# Row 4: 2016/12  (IFRS연결)
$ws.Range("D4").Value = 114510
$ws.Range("E4").Value = 7465
$ws.Range("F4").Value = 7465
$ws.Range("G4").Value = 6426
$ws.Range("H4").Value = 4927
$ws.Range("I4").Value = 4928
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 119891
$ws.Range("L4").Value = 71628
$ws.Range("M4").Value = 48263
$ws.Range("N4").Value = 48261
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 25740
$ws.Range("Q4").Value = 22248
$ws.Range("R4").Value = -14923
$ws.Range("S4").Value = -7073
$ws.Range("T4").Value = 12836
$ws.Range("U4").Value = 9412
$ws.Range("V4").Value = 39792
$ws.Range("W4").Value = 6.52
$ws.Range("X4").Value = 4.3
$ws.Range("Y4").Value = 10.63
$ws.Range("Z4").Value = 4.12
$ws.Range("AA4").Value = 148.41
$ws.Range("AB4").Value = 87.45
$ws.Range("AC4").Value = 1129
$ws.Range("AD4").Value = 10.15
$ws.Range("AE4").Value = 11054
$ws.Range("AF4").Value = 1.04
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 3.06
$ws.Range("AI4").Value = 31.01
$ws.Range("AJ4").Value = 436611361

# Row 5: 2017/12  (IFRS연결)
$ws.Range("D5").Value = 122794
$ws.Range("E5").Value = 8263
$ws.Range("F5").Value = 8263
$ws.Range("G5").Value = 6668
$ws.Range("H5").Value = 5471
$ws.Range("I5").Value = 5472
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 119355
$ws.Range("L5").Value = 67025
$ws.Range("M5").Value = 52330
$ws.Range("N5").Value = 52329
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 25740
$ws.Range("Q5").Value = 22794
$ws.Range("R5").Value = -13583
$ws.Range("S5").Value = -7885
$ws.Range("T5").Value = 11826
$ws.Range("U5").Value = 10969
$ws.Range("V5").Value = 33416
$ws.Range("W5").Value = 6.73
$ws.Range("X5").Value = 4.46
$ws.Range("Y5").Value = 10.88
$ws.Range("Z5").Value = 4.57
$ws.Range("AA5").Value = 128.08
$ws.Range("AB5").Value = 103.24
$ws.Range("AC5").Value = 1253
$ws.Range("AD5").Value = 11.17
$ws.Range("AE5").Value = 11985
$ws.Range("AF5").Value = 1.17
$ws.Range("AG5").Value = 400
$ws.Range("AH5").Value = 2.86
$ws.Range("AI5").Value = 31.91
$ws.Range("AJ5").Value = 436611361

# Row 6: 2018/12  (IFRS연결)
$ws.Range("D6").Value = 121251
$ws.Range("E6").Value = 7309
$ws.Range("F6").Value = 7309
$ws.Range("G6").Value = 6483
$ws.Range("H6").Value = 4816
$ws.Range("I6").Value = 4816
$ws.Range("K6").Value = 139399
$ws.Range("L6").Value = 70866
$ws.Range("M6").Value = 68534
$ws.Range("N6").Value = 68532
$ws.Range("P6").Value = 25740
$ws.Range("Q6").Value = 20694
$ws.Range("R6").Value = -15936
$ws.Range("S6").Value = -5469
$ws.Range("T6").Value = 12210
$ws.Range("U6").Value = 8485
$ws.Range("V6").Value = 29711
$ws.Range("W6").Value = 6.03
$ws.Range("X6").Value = 3.97
$ws.Range("Y6").Value = 7.97
$ws.Range("Z6").Value = 3.72
$ws.Range("AA6").Value = 103.4
$ws.Range("AB6").Value = 166.44
$ws.Range("AC6").Value = 1103
$ws.Range("AD6").Value = 16
$ws.Range("AE6").Value = 15696
$ws.Range("AF6").Value = 1.12
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 2.27
$ws.Range("AI6").Value = 36.26
$ws.Range("AJ6").Value = 436611361

# Row 7: 2019/12(E)  (IFRS연결)
$ws.Range("D7").Value = 127371
$ws.Range("E7").Value = 6525
$ws.Range("G7").Value = 5689
$ws.Range("H7").Value = 4271
$ws.Range("I7").Value = 4281
$ws.Range("K7").Value = 146686
$ws.Range("L7").Value = 75602
$ws.Range("M7").Value = 71084
$ws.Range("N7").Value = 71084
$ws.Range("P7").Value = 25740
$ws.Range("Q7").Value = 25523
$ws.Range("R7").Value = -25495
$ws.Range("S7").Value = -138
$ws.Range("T7").Value = 22073
$ws.Range("U7").Value = 668
$ws.Range("W7").Value = 5.12
$ws.Range("X7").Value = 3.35
$ws.Range("Y7").Value = 6.13
$ws.Range("Z7").Value = 2.99
$ws.Range("AA7").Value = 106.36
$ws.Range("AC7").Value = 980
$ws.Range("AD7").Value = 13.51
$ws.Range("AE7").Value = 16281
$ws.Range("AF7").Value = 0.8100000000000001
$ws.Range("AG7").Value = 411
$ws.Range("AH7").Value = 3.1
$ws.Range("AI7").Value = 41.96

# Row 8: 2020/12(E)  (IFRS연결)
$ws.Range("D8").Value = 134012
$ws.Range("E8").Value = 7405
$ws.Range("G8").Value = 6584
$ws.Range("H8").Value = 4991
$ws.Range("I8").Value = 4982
$ws.Range("K8").Value = 151578
$ws.Range("L8").Value = 77065
$ws.Range("M8").Value = 74512
$ws.Range("N8").Value = 74286
$ws.Range("P8").Value = 25740
$ws.Range("Q8").Value = 25714
$ws.Range("R8").Value = -23280
$ws.Range("S8").Value = -1564
$ws.Range("T8").Value = 20970
$ws.Range("U8").Value = 1839
$ws.Range("W8").Value = 5.53
$ws.Range("X8").Value = 3.72
$ws.Range("Y8").Value = 6.86
$ws.Range("Z8").Value = 3.35
$ws.Range("AA8").Value = 103.43
$ws.Range("AC8").Value = 1141
$ws.Range("AD8").Value = 11.61
$ws.Range("AE8").Value = 17014
$ws.Range("AF8").Value = 0.78
$ws.Range("AG8").Value = 405
$ws.Range("AH8").Value = 3.05
$ws.Range("AI8").Value = 35.45

# Row 9: 2021/12(E)  (IFRS연결)
$ws.Range("D9").Value = 140822
$ws.Range("E9").Value = 8784
$ws.Range("G9").Value = 8078
$ws.Range("H9").Value = 6138
$ws.Range("I9").Value = 6123
$ws.Range("K9").Value = 157319
$ws.Range("L9").Value = 78566
$ws.Range("M9").Value = 78753
$ws.Range("N9").Value = 78528
$ws.Range("P9").Value = 25740
$ws.Range("Q9").Value = 27625
$ws.Range("R9").Value = -22806
$ws.Range("S9").Value = -2478
$ws.Range("T9").Value = 19631
$ws.Range("U9").Value = 3820
$ws.Range("W9").Value = 6.24
$ws.Range("X9").Value = 4.36
$ws.Range("Y9").Value = 8.01
$ws.Range("Z9").Value = 3.97
$ws.Range("AA9").Value = 99.76000000000001
$ws.Range("AC9").Value = 1402
$ws.Range("AD9").Value = 9.449999999999999
$ws.Range("AE9").Value = 17986
$ws.Range("AF9").Value = 0.74
$ws.Range("AG9").Value = 425
$ws.Range("AH9").Value = 3.21
$ws.Range("AI9").Value = 30.31
